# repull data, push all data, mean calculation
# Update column F (dSF) values for several rows on Sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = -8
    4  = -6
    5  = -8
    6  = -10
    8  = -8
    9  = -2
    11 = -1
    14 = -4
    17 = 1
    21 = 1
    26 = 2
    29 = -2
    36 = -6
    40 = 0
    46 = 5
    52 = 4
    53 = -3
    55 = 4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
